# Refactored settings, added punctuation = accept
#
# Applies to the "Tests" worksheet (xl/worksheets/sheet1.xml):
#  - truncates the long remark in I36 (shared string) to drop the trailing
#    "but possibly too slow nl corpus"
#  - adds new experiment rows 39, 40, 42, 43, 45 (rows 38, 41, 44 stay blank)
#  - updates dimension / selection to the new extent

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tests")

# --- Shorten the existing remark text on row 36 ---
$ws.Range("I36").Value = "Basic set-up, but timblserver now and new calculation of cks and skks…"

# --- Row 39: Sonar1 / Standardtest2 ---
$ws.Cells.Item(39, 1).Value = "Sonar1"
$ws.Cells.Item(39, 2).Value = "Standardtest2"
$ws.Cells.Item(39, 3).Value = 4500
$ws.Cells.Item(39, 4).Value = 4500
$ws.Cells.Item(39, 5).Value = "w"
$ws.Cells.Item(39, 6).Value = 15
$ws.Cells.Item(39, 7).Value = 20
$ws.Cells.Item(39, 8).Value = 20
$ws.Cells.Item(39, 9).Value = "Same"

# --- Row 40: Sonar1 / Standardtest2 - low attenuation, data cols first ---
$ws.Cells.Item(40, 1).Value = "Sonar1"
$ws.Cells.Item(40, 2).Value = "Standardtest2"
$ws.Cells.Item(40, 3).Value = 10
$ws.Cells.Item(40, 4).Value = 10
$ws.Cells.Item(40, 5).Value = "w"
$ws.Cells.Item(40, 6).Value = 20
$ws.Cells.Item(40, 7).Value = 24
$ws.Cells.Item(40, 8).Value = 20

# row 41 intentionally left blank (matches source spreadsheet's spacing)

# --- Row 42: Wessel1 / Standardtest2 - switched safety net, data cols first ---
$ws.Cells.Item(42, 1).Value = "Wessel1"
$ws.Cells.Item(42, 2).Value = "Standardtest2"
$ws.Cells.Item(42, 3).Value = 3
$ws.Cells.Item(42, 4).Value = 3
$ws.Cells.Item(42, 5).Value = "w"
$ws.Cells.Item(42, 6).Value = 21
$ws.Cells.Item(42, 7).Value = 23
$ws.Cells.Item(42, 8).Value = 13

# --- Row 43: Wessel1 / 0.1 attenuation ---
$ws.Cells.Item(43, 1).Value = "Wessel1"
$ws.Cells.Item(43, 2).Value = 0.1
$ws.Range("B43").NumberFormat = $ws.Range("B37").NumberFormat
$ws.Cells.Item(43, 3).Value = 3
$ws.Cells.Item(43, 4).Value = 3
$ws.Cells.Item(43, 5).Value = "w"
$ws.Cells.Item(43, 6).Value = 33
$ws.Cells.Item(43, 7).Value = 41
$ws.Cells.Item(43, 8).Value = 2329

# row 44 intentionally left blank (matches source spreadsheet's spacing)

# --- Row 45: Wessel1 / 0.1 attenuation - added punctuation = accept prediction ---
$ws.Cells.Item(45, 1).Value = "Wessel1"
$ws.Cells.Item(45, 2).Value = 0.1
$ws.Range("B45").NumberFormat = $ws.Range("B37").NumberFormat
$ws.Cells.Item(45, 3).Value = 3
$ws.Cells.Item(45, 4).Value = 3
$ws.Cells.Item(45, 5).Value = "w"
$ws.Cells.Item(45, 6).Value = 38
$ws.Cells.Item(45, 7).Value = 44
$ws.Cells.Item(45, 8).Value = 2696

# --- Remarks column (I), filled in the same order the original author typed
#     them so newly-interned shared strings line up with the source file ---
$ws.Cells.Item(42, 9).Value = "Switched safety net to the low-attenuation one"
$ws.Cells.Item(40, 9).Value = "Same, notice low attenuation value now possible"
$ws.Cells.Item(45, 9).Value = "Added punctuation = accept prediction"

# --- Move the active selection to the new last cell, like the source workbook ---
$ws.Activate() | Out-Null
$ws.Range("I45").Select() | Out-Null
